$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.086.51"
$ws.Range("E2").Value = "  -0.61%  "

$ws.Range("D3").Value = "1.805.75"
$ws.Range("E3").Value = "  -0.01%  "

$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.16%  "

$ws.Range("D5").Value = "310.44"
$ws.Range("E5").Value = "  -1.32%  "

$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  -0.02%  "

$ws.Range("D7").Value = "0.5087"
$ws.Range("E7").Value = "  -3.54%  "

$ws.Range("D8").Value = "0.3852"
$ws.Range("E8").Value = "  +0.24%  "

$ws.Range("D9").Value = "0.08422"
$ws.Range("E9").Value = "  +4.77%  "

$ws.Range("D10").Value = "1.097"
$ws.Range("E10").Value = "  -0.65%  "

$ws.Range("D11").Value = "40.80"
$ws.Range("E11").Value = "  -1.79%  "

$ws.Range("D12").Value = "6.382"
$ws.Range("E12").Value = "  +0.40%  "

$ws.Range("D13").Value = "0.9999"
$ws.Range("E13").Value = "  -0.16%  "

$ws.Range("D14").Value = "20.39"
$ws.Range("E14").Value = "  -1.48%  "

$ws.Range("D15").Value = "1.801.88"
$ws.Range("E15").Value = "  -0.09%  "

$ws.Range("D16").Value = "7.291"
$ws.Range("E16").Value = "  -0.89%  "

$ws.Range("D17").Value = "92.25"
$ws.Range("E17").Value = "  -0.70%  "

$ws.Range("D18").Value = "0.00001091"
$ws.Range("E18").Value = "  -0.68%  "

$ws.Range("D19").Value = "0.06586"
$ws.Range("E19").Value = "  -0.33%  "

$ws.Range("D20").Value = "1.001"
$ws.Range("E20").Value = "  +0.02%  "

$ws.Range("D21").Value = "17.31"
$ws.Range("E21").Value = "  -0.77%  "

$ws.Range("D22").Value = "6.007"
$ws.Range("E22").Value = "  +0.27%  "

$ws.Range("D23").Value = "28.085.59"
$ws.Range("E23").Value = "  -0.74%  "

$ws.Range("D24").Value = "11.08"
$ws.Range("E24").Value = "  -1.52%  "

$ws.Range("D25").Value = "2.221"
$ws.Range("E25").Value = "  -0.58%  "

$ws.Range("D26").Value = "158.70"
$ws.Range("E26").Value = "  -0.96%  "

$ws.Range("B27").Value = "LidoDAOToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D27").Value = "2.418"
$ws.Range("E27").Value = "  +0.92%  "

$ws.Range("B28").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C28").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D28").Value = "2.010.24"
$ws.Range("E28").Value = "  -0.01%  "

$ws.Range("D29").Value = "20.37"
$ws.Range("E29").Value = "  -1.03%  "

$ws.Range("D30").Value = "127.62"
$ws.Range("E30").Value = "  +3.44%  "

$ws.Range("D31").Value = "0.1093"
$ws.Range("E31").Value = "  +0.36%  "

$ws.Range("D32").Value = "1.048"
$ws.Range("E32").Value = "  -1.48%  "

$ws.Range("D33").Value = "3.649"
$ws.Range("E33").Value = "  -0.53%  "

$ws.Range("D34").Value = "5.570"
$ws.Range("E34").Value = "  -0.29%  "

$ws.Range("D35").Value = "0.06970"
$ws.Range("E35").Value = "  -4.29%  "

$ws.Range("D36").Value = "9.139"
$ws.Range("E36").Value = "  +2.97%  "

$ws.Range("D37").Value = "0.02336"
$ws.Range("E37").Value = "  +0.67%  "

$ws.Range("D38").Value = "0.2174"
$ws.Range("E38").Value = "  +0.05%  "

$ws.Range("D39").Value = "5.017"
$ws.Range("E39").Value = "  -1.96%  "

$ws.Range("D40").Value = "11.44"
$ws.Range("E40").Value = "  -7.28%  "

$ws.Range("D41").Value = "0.6130"
$ws.Range("E41").Value = "  -1.70%  "

$ws.Range("E42").Value = "  +0.10%  "

$ws.Range("D43").Value = "1.152"
$ws.Range("E43").Value = "  -1.37%  "

$ws.Range("D44").Value = "13.25"
$ws.Range("E44").Value = "  +0.04%  "

$ws.Range("D45").Value = "1.297"
$ws.Range("E45").Value = "  -5.31%  "

$ws.Range("D46").Value = "0.5907"
$ws.Range("E46").Value = "  -2.03%  "

$ws.Range("D47").Value = "3.711"
$ws.Range("E47").Value = "  -1.53%  "

$ws.Range("D48").Value = "125.76"
$ws.Range("E48").Value = "  -1.03%  "

$ws.Range("D49").Value = "1.935"
$ws.Range("E49").Value = "  -0.12%  "

$ws.Range("D50").Value = "1.187"
$ws.Range("E50").Value = "  -2.29%  "

$ws.Range("D51").Value = "0.06735"
$ws.Range("E51").Value = "  -1.57%  "
